# edit.ps1 - Applies the Journal.docx OOXML diff via Word COM-interop (PowerShell style)
#
# Summary of required changes:
#  1. Several paragraphs have their runs (previously split apart by
#     <w:proofErr .../> spell/grammar-check markers) merged back into a
#     single contiguous run, with the proofErr markers removed.
#  2. The final paragraph ("Will work on the progress report this weekend")
#     loses its <w:proofErr w:type="gramStart/gramEnd"/> wrapper and the
#     "_GoBack" bookmark moves off of it.
#  3. A new "October 11 – October 30" section (one Heading1 paragraph plus
#     three list-item paragraphs) is appended at the end of the document,
#     with the "_GoBack" bookmark now sitting on the very last paragraph.

$d = $word.ActiveDocument

function Wrap-Pkg([string]$bodyXml) {
    # Word's Range.InsertXML wants either a full WordOpenXML "package" or a
    # fragment shaped like one; bare <w:r>/<w:p> fragments are not reliably
    # accepted, so every insert below is wrapped the same way Range.WordOpenXML
    # itself would emit.
    return '<?xml version="1.0" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-RangeWithRun([object]$range, [string]$text, [bool]$preserve) {
    # Collapses the given Range down to a single run containing $text
    # (stripping out any proofErr markers / run-splits that previously
    # lived inside that range).
    if ($preserve) {
        $t = '<w:t xml:space="preserve">' + $text + '</w:t>'
    } else {
        $t = '<w:t>' + $text + '</w:t>'
    }
    $xml = Wrap-Pkg ('<w:p><w:r>' + $t + '</w:r></w:p>')
    $range.InsertXML($xml)
}

function Find-RangeFor([object]$searchRange, [string]$needle) {
    # Locates $needle inside $searchRange and returns a new collapsed-to-match
    # Range (mirrors what Find.Execute does to its own Range in place).
    $hit = $d.Range($searchRange.Start, $searchRange.End)
    [void]$hit.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $hit
}

# ---------------------------------------------------------------------
# 1) "These fragments needed to be managed with a FragmentPagerAdapter"
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "These fragments needed to be managed with a "
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
Replace-RangeWithRun $full "These fragments needed to be managed with a FragmentPagerAdapter" $false

# ---------------------------------------------------------------------
# 2) "The FragmentPagerAdapter deals with how the fragments are stored in memory"
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "FragmentPagerAdapter deals"
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
Replace-RangeWithRun $full "The FragmentPagerAdapter deals with how the fragments are stored in memory" $false

# ---------------------------------------------------------------------
# 3) "Each fragment requires a layout and needs to be stored in a container (viewpager)"
#    -- only the " in a container (viewpager)" tail gets collapsed; the
#    leading two runs ("...and need" / "s to be stored") are untouched.
# ---------------------------------------------------------------------
$tail = Find-RangeFor $d.Content " in a container ("
$paraEnd = $tail.Paragraphs(1).Range.End
$full = $d.Range($tail.Start, $paraEnd - 1)
Replace-RangeWithRun $full " in a container (viewpager)" $true

# ---------------------------------------------------------------------
# 4) "The viewpager takes care of switching to the correct layout on correct screen"
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "viewpager takes care"
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
Replace-RangeWithRun $full "The viewpager takes care of switching to the correct layout on correct screen" $false

# ---------------------------------------------------------------------
# 5) I have also got the hang of ... setting “onClickListeners” for different widgets
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "I have also got the hang of"
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
$text5 = "I have also got the hang of the common xml tags such as ID, width, height etc. and setting “onClickListeners” for different widgets"
Replace-RangeWithRun $full $text5 $false

# ---------------------------------------------------------------------
# 6) "Talked to Mr. Grondin about learning angular, js, node and mongoDB"
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "Talked to Mr."
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
Replace-RangeWithRun $full "Talked to Mr. Grondin about learning angular, js, node and mongoDB" $false

# ---------------------------------------------------------------------
# 7) "Will work on the progress report this weekend" -- merge text, drop
#    the gramStart/gramEnd proofErr wrapper, and strip the _GoBack
#    bookmark from this paragraph (it gets re-created on the new last
#    paragraph below).
# ---------------------------------------------------------------------
$r = Find-RangeFor $d.Content "Will work on the progress report this "
$p = $r.Paragraphs(1).Range
$full = $d.Range($p.Start, $p.End - 1)
Replace-RangeWithRun $full "Will work on the progress report this weekend" $false

# ---------------------------------------------------------------------
# 8) Append the new "October 11 – October 30" section at the end of the
#    document, with the _GoBack bookmark on the final paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newSection = '<w:p/>' + `
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>October 11 – October 30</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Hoping to finish the tutorial series between now and the next report</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Learned about interfaces and dialogs in android studio</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Can update email and username in the database</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$endRange.InsertXML((Wrap-Pkg $newSection))

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
